$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a "text-like" value into a cell while preserving the cell's
# existing number format / style (copied beforehand from a template cell).
# Many of the values in this sheet (counts like "1", "0", prices like
# "123.00") must be stored as shared-string TEXT, not as native numbers,
# to match the original workbook's layout. Temporarily switching the
# format to "@" (text) forces Excel to store the value as text; we then
# restore the original formatting by pasting the format (only) back from
# a same-column template cell, without touching the value we just set.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddr, $templateAddr, $value) {
    $ws.Range($rangeAddr).NumberFormat = "@"
    $ws.Range($rangeAddr).Value = $value
    $ws.Range($templateAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 1. Insert 3 new blank rows before the totals row (current row 16), so the
#    existing totals (row16->19) and footer (row17->20) rows shift down.
# ---------------------------------------------------------------------------
$ws.Range("A16:A18").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Give the new rows 16-18 the same column formatting as row 15 (a regular
#    product row), then fix up row heights and merged cells to match.
# ---------------------------------------------------------------------------
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)
$ws.Range("A17:Q17").PasteSpecial(-4122)
$ws.Range("A18:Q18").PasteSpecial(-4122)

$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75

foreach ($r in 16,17,18) {
    $ws.Range("A$r" + ":B$r").Merge()
    $ws.Range("C$r" + ":G$r").Merge()
    $ws.Range("H$r" + ":K$r").Merge()
    $ws.Range("L$r" + ":M$r").Merge()
    $ws.Range("N$r" + ":O$r").Merge()
}

# ---------------------------------------------------------------------------
# 3. Row 15 used to be the last product ("شمع حريمي"); it now becomes the
#    new product "TORSERETIC 20MG 30 TABS.".
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "TORSERETIC 20MG 30 TABS."
Set-TextValue "H15" "H15" "1:2"
Set-TextValue "L15" "L15" "1"
Set-TextValue "N15" "N15" "123.00"
Set-TextValue "P15" "P15" "123.0000"
Set-TextValue "Q15" "Q15" "1:0"

# ---------------------------------------------------------------------------
# 4. Row 16: new product "حبايه".
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "حبايه"
Set-TextValue "H16" "H16" "0:0"
Set-TextValue "L16" "L16" "0"
Set-TextValue "N16" "N16" "3.00"
Set-TextValue "P16" "P16" "6.0000"
Set-TextValue "Q16" "Q16" "2:0"

# ---------------------------------------------------------------------------
# 5. Row 17: product "شمع حريمي" (the product that used to live in row 15).
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "شمع حريمي"
Set-TextValue "H17" "H17" "5:0"
Set-TextValue "L17" "L17" "0"
Set-TextValue "N17" "N17" "50.00"
Set-TextValue "P17" "P17" "50.0000"
Set-TextValue "Q17" "Q17" "1:0"

# ---------------------------------------------------------------------------
# 6. Row 18: new product "مجموعه برد".
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = 12
$ws.Range("C18").Value = "مجموعه برد"
Set-TextValue "H18" "H18" "0:0"
Set-TextValue "L18" "L18" "0"
Set-TextValue "N18" "N18" "8.00"
Set-TextValue "P18" "P18" "8.0000"
Set-TextValue "Q18" "Q18" "1:0"

# ---------------------------------------------------------------------------
# 7. The totals row (previously row 16, now row 19) needs an updated total.
# ---------------------------------------------------------------------------
$ws.Range("P19").Value = 509.87
